{"js": "// The document contains a \"HIVE TEAMS\" roster. This edit translates a\n// handful of English role labels/snippets in the \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\" (Outreach)\n// team section into Russian, matching the author's new translations.\n//\n// Changes (document order):\n//   1. \"Editor of \"\" (run, before the \"Your Week in SmartCash\" hyperlink)\n//          -> \"\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"\"\n//   2. Paragraph \"Outreach Support\" right after \"Semptly\"\n//          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n//   3. Paragraph \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\" right after \"Carlos Santiago\"\n//          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\"\n//   4. Paragraph \"Outreach Support\" right after \"illumin8\"\n//          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n//\n// NOTE: the later \"HIVE TEAM: OUTREACH 2\" section also contains two\n// \"Outreach Support\" paragraphs (for \"Eiky\" and \"Filipe Boldo\") that must\n// stay untouched, so we do not blanket-replace every match \u2014 we target\n// specific paragraphs located relative to their neighboring name headings.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Helper: get the plain text of a paragraph.\nasync function textOf(p) {\n  p.load(\"text\");\n  await context.sync();\n  return p.text;\n}\n\n// Helper: replace the *first* run matching `oldText` inside paragraph `p`\n// with `newText`, preserving the run's own formatting.\nasync function replaceInParagraph(p, oldText, newText) {\n  const found = p.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  if (found.items.length === 0) {\n    throw new Error(`Text \"${oldText}\" not found in target paragraph.`);\n  }\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) \"Editor of \"\" -> \"\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"\" (keeps the hyperlink run right after it).\nfor (const p of paragraphs.items) {\n  const t = await textOf(p);\n  if (t.indexOf('Editor of \"') === 0) {\n    await replaceInParagraph(p, 'Editor of \"', '\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"');\n    break;\n  }\n}\n\n// Re-load paragraph texts/order after the edit above (paragraph count is\n// unchanged, but refresh to be safe before locating the remaining targets).\nparagraphs.load(\"items\");\nawait context.sync();\nconst names = [];\nfor (const p of paragraphs.items) {\n  names.push(await textOf(p));\n}\n\n// 2) The \"Outreach Support\" paragraph immediately following \"Semptly\".\nconst semptlyIdx = names.indexOf(\"Semptly\");\nif (semptlyIdx === -1 || names[semptlyIdx + 1] !== \"Outreach Support\") {\n  throw new Error(\"Could not locate the 'Semptly' / 'Outreach Support' pair.\");\n}\nawait replaceInParagraph(\n  paragraphs.items[semptlyIdx + 1],\n  \"Outreach Support\",\n  \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n);\n\n// 3) The \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\" paragraph immediately following\n//    \"Carlos Santiago\".\nconst carlosIdx = names.indexOf(\"Carlos Santiago\");\nif (carlosIdx === -1 || names[carlosIdx + 1] !== \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\") {\n  throw new Error(\"Could not locate the 'Carlos Santiago' / '\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433' pair.\");\n}\nawait replaceInParagraph(\n  paragraphs.items[carlosIdx + 1],\n  \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\",\n  \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\"\n);\n\n// 4) The \"Outreach Support\" paragraph immediately following \"illumin8\"\n//    (there is another, unrelated \"Outreach Support\" pair later in the\n//    document under \"HIVE TEAM: OUTREACH 2\" that must stay untouched).\nconst illumin8Idx = names.indexOf(\"illumin8\");\nif (illumin8Idx === -1 || names[illumin8Idx + 1] !== \"Outreach Support\") {\n  throw new Error(\"Could not locate the 'illumin8' / 'Outreach Support' pair.\");\n}\nawait replaceInParagraph(\n  paragraphs.items[illumin8Idx + 1],\n  \"Outreach Support\",\n  \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n);\n", "ps1": "# The document contains a \"HIVE TEAMS\" roster. This edit translates a\n# handful of English role labels/snippets in the \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\" (Outreach)\n# team section into Russian, matching the author's new translations.\n#\n# Changes (document order):\n#   1. \"Editor of \"\"\" (run, right before the \"Your Week in SmartCash\" hyperlink)\n#          -> \"\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"\"\"\n#   2. Paragraph \"Outreach Support\" right after \"Semptly\"\n#          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n#   3. Paragraph \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\" right after \"Carlos Santiago\"\n#          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\"\n#   4. Paragraph \"Outreach Support\" right after \"illumin8\"\n#          -> \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n#\n# NOTE: a later \"HIVE TEAM: OUTREACH 2\" section also contains two\n# \"Outreach Support\" paragraphs (for \"Eiky\" and \"Filipe Boldo\") that must\n# stay untouched, so we do not blanket Find/Replace every match across the\n# whole document \u2014 we target specific paragraphs located relative to their\n# neighboring name headings. We also avoid Find.Execute's Replacement\n# parameter (it smart-quotes straight quotes via AutoFormat/AutoCorrect)\n# and instead assign Range.Text directly.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaTexts($doc) {\n    $texts = New-Object System.Collections.ArrayList\n    foreach ($p in $doc.Paragraphs) {\n        [void]$texts.Add($p.Range.Text)\n    }\n    return $texts\n}\n\n# 1) \"Editor of \"\"\" -> \"\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"\"\" (keep the hyperlink run right after it).\n$found = $false\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range.Duplicate\n    if ($rng.Find.Execute('Editor of \"', $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)) {\n        $rng.Text = '\u0420\u0435\u0434\u0430\u043a\u0442\u043e\u0440 \"'\n        $found = $true\n        break\n    }\n}\nif (-not $found) {\n    throw \"Could not locate the 'Editor of \"\"' run.\"\n}\n\n# 2) The \"Outreach Support\" paragraph immediately following \"Semptly\".\n$paras = $d.Paragraphs\n$n = $paras.Count\n$semptlyIdx = -1\nfor ($i = 1; $i -le $n; $i++) {\n    if ($paras.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq \"Semptly\") {\n        $semptlyIdx = $i\n        break\n    }\n}\nif ($semptlyIdx -eq -1) {\n    throw \"Could not locate the 'Semptly' paragraph.\"\n}\n$nextPara = $paras.Item($semptlyIdx + 1)\n$nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)\nif ($nextText -ne \"Outreach Support\") {\n    throw \"Expected 'Outreach Support' after 'Semptly', found '$nextText'.\"\n}\n$rng = $nextPara.Range.Duplicate\n$rng.MoveEnd(1, -1) | Out-Null\n$rng.Text = \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n\n# 3) The \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\" paragraph immediately following\n#    \"Carlos Santiago\".\n$paras = $d.Paragraphs\n$n = $paras.Count\n$carlosIdx = -1\nfor ($i = 1; $i -le $n; $i++) {\n    if ($paras.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq \"Carlos Santiago\") {\n        $carlosIdx = $i\n        break\n    }\n}\nif ($carlosIdx -eq -1) {\n    throw \"Could not locate the 'Carlos Santiago' paragraph.\"\n}\n$nextPara = $paras.Item($carlosIdx + 1)\n$nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)\nif ($nextText -ne \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\") {\n    throw \"Expected '\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433' after 'Carlos Santiago', found '$nextText'.\"\n}\n$rng = $nextPara.Range.Duplicate\n$rng.MoveEnd(1, -1) | Out-Null\n$rng.Text = \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435\"\n\n# 4) The \"Outreach Support\" paragraph immediately following \"illumin8\"\n#    (there is another, unrelated \"Outreach Support\" pair later in the\n#    document under \"HIVE TEAM: OUTREACH 2\" that must stay untouched).\n$paras = $d.Paragraphs\n$n = $paras.Count\n$illumin8Idx = -1\nfor ($i = 1; $i -le $n; $i++) {\n    if ($paras.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq \"illumin8\") {\n        $illumin8Idx = $i\n        break\n    }\n}\nif ($illumin8Idx -eq -1) {\n    throw \"Could not locate the 'illumin8' paragraph.\"\n}\n$nextPara = $paras.Item($illumin8Idx + 1)\n$nextText = $nextPara.Range.Text.TrimEnd([char]13, [char]7)\nif ($nextText -ne \"Outreach Support\") {\n    throw \"Expected 'Outreach Support' after 'illumin8', found '$nextText'.\"\n}\n$rng = $nextPara.Range.Duplicate\n$rng.MoveEnd(1, -1) | Out-Null\n$rng.Text = \"\u041f\u0440\u043e\u0434\u0432\u0438\u0436\u0435\u043d\u0438\u0435 \u0438 \u043c\u0430\u0440\u043a\u0435\u0442\u0438\u043d\u0433\"\n"}
